# Auto-generated COM-interop edit script
# Applies the checkpoint edit to genx_signals.xlsx (3 sheets)

$wb = $excel.ActiveWorkbook

$wsActive = $wb.Worksheets.Item("Active Signals")

# Row 2
$wsActive.Cells.Item(2,1).Value = '2025-07-28 21:03'
$wsActive.Cells.Item(2,2).Value = 'EURUSD'
$wsActive.Cells.Item(2,3).Value = 'BUY'
$wsActive.Cells.Item(2,4).Value = 1.10695
$wsActive.Cells.Item(2,5).Value = 1.10383
$wsActive.Cells.Item(2,6).Value = 1.11219
$wsActive.Cells.Item(2,7).Value = 0.07
$wsActive.Cells.Item(2,8).NumberFormat = "@"
$wsActive.Cells.Item(2,8).Value = '81.0%'
$wsActive.Cells.Item(2,9).Value = 1.68
$wsActive.Cells.Item(2,10).Value = 'Active'

# Row 3
$wsActive.Cells.Item(3,1).Value = '2025-07-28 21:16'
$wsActive.Cells.Item(3,2).Value = 'USDCHF'
$wsActive.Cells.Item(3,3).Value = 'SELL'
$wsActive.Cells.Item(3,4).Value = 0.87879
$wsActive.Cells.Item(3,5).Value = 0.88296
$wsActive.Cells.Item(3,6).Value = 0.87072
$wsActive.Cells.Item(3,7).Value = 0.04
$wsActive.Cells.Item(3,8).NumberFormat = "@"
$wsActive.Cells.Item(3,8).Value = '79.0%'
$wsActive.Cells.Item(3,9).Value = 1.93
$wsActive.Cells.Item(3,10).Value = 'Active'

# Row 4
$wsActive.Cells.Item(4,1).Value = '2025-07-28 20:40'
$wsActive.Cells.Item(4,2).Value = 'USDCHF'
$wsActive.Cells.Item(4,3).Value = 'BUY'
$wsActive.Cells.Item(4,4).Value = 0.87926
$wsActive.Cells.Item(4,5).Value = 0.87699
$wsActive.Cells.Item(4,6).Value = 0.8847
$wsActive.Cells.Item(4,7).Value = 0.04
$wsActive.Cells.Item(4,8).NumberFormat = "@"
$wsActive.Cells.Item(4,8).Value = '80.0%'
$wsActive.Cells.Item(4,9).Value = 2.4
$wsActive.Cells.Item(4,10).Value = 'Active'

# Swap the SELL/BUY highlight fill between C3 and C4 (style-index swap,
# reuses the workbook's existing two fills rather than creating new ones)
$wsActive.Range("C3").Copy()
$wsActive.Range("Z1").PasteSpecial(-4122)
$wsActive.Range("C4").Copy()
$wsActive.Range("C3").PasteSpecial(-4122)
$wsActive.Range("Z1").Copy()
$wsActive.Range("C4").PasteSpecial(-4122)
$wsActive.Range("Z1").Clear()

# Drop the two rows that are no longer active signals
$wsActive.Range("A5:A6").EntireRow.Delete()

$wsDash = $wb.Worksheets.Item("Summary Dashboard")
$wsDash.Range("B4").Value = 3.0
$wsDash.Range("B5").Value = 5.0
$wsDash.Range("B6").Value = 10.0
$wsDash.Range("B7").NumberFormat = "@"
$wsDash.Range("B7").Value = '85.8%'
$wsDash.Range("B8").NumberFormat = "@"
$wsDash.Range("B8").Value = '2.14'
$wsDash.Range("B9").NumberFormat = "@"
$wsDash.Range("B9").Value = '2025-07-28 20:56:28'

$wsHist = $wb.Worksheets.Item("Signal History")

# Row 2
$wsHist.Cells.Item(2,1).Value = '2025-07-28 21:03'
$wsHist.Cells.Item(2,2).Value = 'EURUSD'
$wsHist.Cells.Item(2,3).Value = 'BUY'
$wsHist.Cells.Item(2,4).Value = 1.10695
$wsHist.Cells.Item(2,5).Value = 1.10383
$wsHist.Cells.Item(2,6).Value = 1.11219
$wsHist.Cells.Item(2,7).Value = 0.07
$wsHist.Cells.Item(2,8).Value = 0.81
$wsHist.Cells.Item(2,9).Value = 1.68
$wsHist.Cells.Item(2,10).Value = 'Active'

# Row 3
$wsHist.Cells.Item(3,1).Value = '2025-07-28 20:53'
$wsHist.Cells.Item(3,2).Value = 'USDCHF'
$wsHist.Cells.Item(3,3).Value = 'SELL'
$wsHist.Cells.Item(3,4).Value = 0.88184
$wsHist.Cells.Item(3,5).Value = 0.8858
$wsHist.Cells.Item(3,6).Value = 0.87689
$wsHist.Cells.Item(3,7).Value = 0.1
$wsHist.Cells.Item(3,8).Value = 0.89
$wsHist.Cells.Item(3,9).Value = 1.25
$wsHist.Cells.Item(3,10).Value = 'Pending'

# Row 4
$wsHist.Cells.Item(4,1).Value = '2025-07-28 20:32'
$wsHist.Cells.Item(4,2).Value = 'EURUSD'
$wsHist.Cells.Item(4,3).Value = 'BUY'
$wsHist.Cells.Item(4,4).Value = 1.10663
$wsHist.Cells.Item(4,5).Value = 1.10299
$wsHist.Cells.Item(4,6).Value = 1.11359
$wsHist.Cells.Item(4,7).Value = 0.08
$wsHist.Cells.Item(4,8).Value = 0.93
$wsHist.Cells.Item(4,9).Value = 1.91
$wsHist.Cells.Item(4,10).Value = 'Pending'

# Row 5
$wsHist.Cells.Item(5,1).Value = '2025-07-28 21:13'
$wsHist.Cells.Item(5,2).Value = 'USDJPY'
$wsHist.Cells.Item(5,3).Value = 'BUY'
$wsHist.Cells.Item(5,4).Value = 149.23577
$wsHist.Cells.Item(5,5).Value = 149.00513
$wsHist.Cells.Item(5,6).Value = 149.6558
$wsHist.Cells.Item(5,7).Value = 0.03
$wsHist.Cells.Item(5,8).Value = 0.81
$wsHist.Cells.Item(5,9).Value = 1.82
$wsHist.Cells.Item(5,10).Value = 'Active'

# Row 6
$wsHist.Cells.Item(6,1).Value = '2025-07-28 20:30'
$wsHist.Cells.Item(6,2).Value = 'XAUAUD'
$wsHist.Cells.Item(6,3).Value = 'SELL'
$wsHist.Cells.Item(6,4).Value = 4065.77924
$wsHist.Cells.Item(6,5).Value = 4065.78308
$wsHist.Cells.Item(6,6).Value = 4065.77082
$wsHist.Cells.Item(6,7).Value = 0.04
$wsHist.Cells.Item(6,8).Value = 0.83
$wsHist.Cells.Item(6,9).Value = 2.19
$wsHist.Cells.Item(6,10).Value = 'Filled'

# Row 7
$wsHist.Cells.Item(7,1).Value = '2025-07-28 21:03'
$wsHist.Cells.Item(7,2).Value = 'XAUCHF'
$wsHist.Cells.Item(7,3).Value = 'SELL'
$wsHist.Cells.Item(7,4).Value = 2344.45541
$wsHist.Cells.Item(7,5).Value = 2344.45954
$wsHist.Cells.Item(7,6).Value = 2344.44827
$wsHist.Cells.Item(7,7).Value = 0.04
$wsHist.Cells.Item(7,8).Value = 0.76
$wsHist.Cells.Item(7,9).Value = 1.73
$wsHist.Cells.Item(7,10).Value = 'Pending'

# Row 8
$wsHist.Cells.Item(8,1).Value = '2025-07-28 21:12'
$wsHist.Cells.Item(8,2).Value = 'EURUSD'
$wsHist.Cells.Item(8,3).Value = 'BUY'
$wsHist.Cells.Item(8,4).Value = 1.10121
$wsHist.Cells.Item(8,5).Value = 1.10352
$wsHist.Cells.Item(8,6).Value = 1.09431
$wsHist.Cells.Item(8,7).Value = 0.02
$wsHist.Cells.Item(8,8).Value = 0.85
$wsHist.Cells.Item(8,9).Value = 2.98
$wsHist.Cells.Item(8,10).Value = 'Filled'

# Row 9
$wsHist.Cells.Item(9,1).Value = '2025-07-28 21:16'
$wsHist.Cells.Item(9,2).Value = 'USDCHF'
$wsHist.Cells.Item(9,3).Value = 'SELL'
$wsHist.Cells.Item(9,4).Value = 0.87879
$wsHist.Cells.Item(9,5).Value = 0.88296
$wsHist.Cells.Item(9,6).Value = 0.87072
$wsHist.Cells.Item(9,7).Value = 0.04
$wsHist.Cells.Item(9,8).Value = 0.79
$wsHist.Cells.Item(9,9).Value = 1.93
$wsHist.Cells.Item(9,10).Value = 'Active'

# Row 10
$wsHist.Cells.Item(10,1).Value = '2025-07-28 21:25'
$wsHist.Cells.Item(10,2).Value = 'XAUUSD'
$wsHist.Cells.Item(10,3).Value = 'BUY'
$wsHist.Cells.Item(10,4).Value = 2640.34649
$wsHist.Cells.Item(10,5).Value = 2640.34899
$wsHist.Cells.Item(10,6).Value = 2640.33792
$wsHist.Cells.Item(10,7).Value = 0.06
$wsHist.Cells.Item(10,8).Value = 0.86
$wsHist.Cells.Item(10,9).Value = 3.41
$wsHist.Cells.Item(10,10).Value = 'Filled'

# Row 11
$wsHist.Cells.Item(11,1).Value = '2025-07-28 20:40'
$wsHist.Cells.Item(11,2).Value = 'USDCHF'
$wsHist.Cells.Item(11,3).Value = 'BUY'
$wsHist.Cells.Item(11,4).Value = 0.87926
$wsHist.Cells.Item(11,5).Value = 0.87699
$wsHist.Cells.Item(11,6).Value = 0.8847
$wsHist.Cells.Item(11,7).Value = 0.04
$wsHist.Cells.Item(11,8).Value = 0.8
$wsHist.Cells.Item(11,9).Value = 2.4
$wsHist.Cells.Item(11,10).Value = 'Active'

# Row 12
$wsHist.Cells.Item(12,1).Value = '2025-07-28 20:28'
$wsHist.Cells.Item(12,2).Value = 'AUDUSD'
$wsHist.Cells.Item(12,3).Value = 'SELL'
$wsHist.Cells.Item(12,4).Value = 0.65717
$wsHist.Cells.Item(12,5).Value = 0.66009
$wsHist.Cells.Item(12,6).Value = 0.6475
$wsHist.Cells.Item(12,7).Value = 0.02
$wsHist.Cells.Item(12,8).Value = 0.87
$wsHist.Cells.Item(12,9).Value = 3.31
$wsHist.Cells.Item(12,10).Value = 'Filled'

# Row 13
$wsHist.Cells.Item(13,1).Value = '2025-07-28 21:26'
$wsHist.Cells.Item(13,2).Value = 'EURUSD'
$wsHist.Cells.Item(13,3).Value = 'SELL'
$wsHist.Cells.Item(13,4).Value = 1.10421
$wsHist.Cells.Item(13,5).Value = 1.10631
$wsHist.Cells.Item(13,6).Value = 1.09805
$wsHist.Cells.Item(13,7).Value = 0.03
$wsHist.Cells.Item(13,8).Value = 0.89
$wsHist.Cells.Item(13,9).Value = 2.93
$wsHist.Cells.Item(13,10).Value = 'Filled'

# Row 14
$wsHist.Cells.Item(14,1).Value = '2025-07-28 21:09'
$wsHist.Cells.Item(14,2).Value = 'USDCAD'
$wsHist.Cells.Item(14,3).Value = 'SELL'
$wsHist.Cells.Item(14,4).Value = 1.36718
$wsHist.Cells.Item(14,5).Value = 1.3706
$wsHist.Cells.Item(14,6).Value = 1.36273
$wsHist.Cells.Item(14,7).Value = 0.04
$wsHist.Cells.Item(14,8).Value = 0.94
$wsHist.Cells.Item(14,9).Value = 1.3
$wsHist.Cells.Item(14,10).Value = 'Pending'

# Row 15
$wsHist.Cells.Item(15,1).Value = '2025-07-28 21:13'
$wsHist.Cells.Item(15,2).Value = 'NZDUSD'
$wsHist.Cells.Item(15,3).Value = 'SELL'
$wsHist.Cells.Item(15,4).Value = 0.58909
$wsHist.Cells.Item(15,5).Value = 0.5918
$wsHist.Cells.Item(15,6).Value = 0.5845
$wsHist.Cells.Item(15,7).Value = 0.04
$wsHist.Cells.Item(15,8).Value = 0.91
$wsHist.Cells.Item(15,9).Value = 1.69
$wsHist.Cells.Item(15,10).Value = 'Pending'

# Row 16
$wsHist.Cells.Item(16,1).Value = '2025-07-28 20:41'
$wsHist.Cells.Item(16,2).Value = 'GBPUSD'
$wsHist.Cells.Item(16,3).Value = 'BUY'
$wsHist.Cells.Item(16,4).Value = 1.27296
$wsHist.Cells.Item(16,5).Value = 1.26797
$wsHist.Cells.Item(16,6).Value = 1.28042
$wsHist.Cells.Item(16,7).Value = 0.04
$wsHist.Cells.Item(16,8).Value = 0.93
$wsHist.Cells.Item(16,9).Value = 1.5
$wsHist.Cells.Item(16,10).Value = 'Pending'

